# daily auto push: 2026-01-18 18:40 UTC
#
# Two new daily-ranking rows ("2026/01/18" / "日" and "2026/01/19" / "月")
# are inserted right after the existing 2026/01/18 block (old row 676),
# pushing every subsequent row down by two. That is modeled here as a
# plain row insert at row 677 (done twice), followed by filling in the
# values for the two freshly inserted rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the old row 677 (2026/12/29 block), shifting
# everything from there on down to rows 679.. / dimension A1:D720.
$ws.Rows.Item(677).Insert()
$ws.Rows.Item(677).Insert()

# New row 677: 2026/01/18, 日, 22, 154
# (leading apostrophe keeps the date-like text literal, same as the other
# date cells in this column, instead of Excel auto-converting it to a date)
$ws.Range("A677").Value = "'2026/01/18"
$ws.Range("B677").Value = "日"
$ws.Range("C677").Value = 22
$ws.Range("D677").Value = 154

# New row 678: 2026/01/19, 月, 1, 169
$ws.Range("A678").Value = "'2026/01/19"
$ws.Range("B678").Value = "月"
$ws.Range("C678").Value = 1
$ws.Range("D678").Value = 169
